$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.834.58"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "2.693.24"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("E9").Value = "  +5.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.401"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "30.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.45%  "
$ws.Range("E14").Value = "  +6.19%  "
$ws.Range("D15").Value = "3.180.18"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").Value = "65.694.06"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "2.700.46"
$ws.Range("E17").Value = "  +8.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "358.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("E21").Value = "  +3.54%  "
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.28%  "
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000106"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.64%  "
$ws.Range("E27").Value = "  +3.07%  "
$ws.Range("E28").Value = "  +5.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.75%  "
$ws.Range("E30").Value = "  +5.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "540.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.89%  "
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.433"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "171.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.14%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  +2.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0614"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("E47").Value = "  +4.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0265"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.655"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0992"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.11%  "
